# TC02_Canine_Filter_StudyType-Transcriptomics.xlsx
# "10 icdc scripts for jenkins"
#
# The FilesTab Cypher query (cell B4 on the "startup" sheet) is rewritten to
# drop the `File Type` and `Breed` columns from its RETURN clause, and the
# sheet's active selection moves from D4 to B4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_type IN ['Transcriptomics']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS ``File Name``, 
        coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesQuery

# Reflect the new active cell/selection (was D4).
[void]$ws.Range("B4").Select()
